# Generate Report for Handback
# Updates the localization-status workbook: the handback has completed and
# is in sync with en-US, so the status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-language handback timestamps
# are refreshed, the stale "version not latest" error is cleared, and the
# columns that now hold the longer status text / cleared error text are
# resized.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.166667
$overview.Columns.Item(6).ColumnWidth = 29.166667

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("K2").Value = "2016-09-04 08:53:40"
$zhcn.Range("P2").Value = ""
$zhcn.Columns.Item(3).ColumnWidth = 29.166667
$zhcn.Columns.Item(16).ColumnWidth = 12.833333

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $statusText
$dede.Range("K2").Value = "2016-09-04 08:53:47"
$dede.Range("P2").Value = ""
$dede.Columns.Item(3).ColumnWidth = 29.166667
$dede.Columns.Item(16).ColumnWidth = 12.833333
